$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the old "Avaliação:" row (row 12) to make room for the
# new "Docentes responsáveis:" section (one header row + two lecturer rows).
$ws.Range("A12:A14").EntireRow.Insert()

# Row 12: section header (column A only, bold style already inherited from Insert)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13: first lecturer, duplicated into column B (normal) and column C (red/changed)
$ws.Range("B13").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C13").Value = "5840726 - Cristina Bormio Nunes"

# Row 14: second lecturer, duplicated into column B and column C
$ws.Range("B14").Value = "1341653 - Maria José Ramos Sandim"
$ws.Range("C14").Value = "1341653 - Maria José Ramos Sandim"

# The row-insert operation auto-fills the whole row width with empty, styled
# cells (copied down from the row above). Remove the cells that should stay
# absent so only the cells with actual content remain, matching the target.
$ws.Range("B12:C12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
